$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-08-29 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-08-30 Saturday", 2) | Out-Null
$d.Content.Find.Execute("87-54=", $true, $false, $false, $false, $false, $true, 1, $false, "34-9=", 2) | Out-Null
$d.Content.Find.Execute("35-26=", $true, $false, $false, $false, $false, $true, 1, $false, "40-32=", 2) | Out-Null
$d.Content.Find.Execute("31+46=", $true, $false, $false, $false, $false, $true, 1, $false, "88-23=", 2) | Out-Null
$d.Content.Find.Execute("2+85=", $true, $false, $false, $false, $false, $true, 1, $false, "19+25=", 2) | Out-Null
$d.Content.Find.Execute("17+31=", $true, $false, $false, $false, $false, $true, 1, $false, "46+1=", 2) | Out-Null
$d.Content.Find.Execute("44-0=", $true, $false, $false, $false, $false, $true, 1, $false, "38+1=", 2) | Out-Null
$d.Content.Find.Execute("17+46=", $true, $false, $false, $false, $false, $true, 1, $false, "64+11=", 2) | Out-Null
$d.Content.Find.Execute("21+16=", $true, $false, $false, $false, $false, $true, 1, $false, "31+13=", 2) | Out-Null
$d.Content.Find.Execute("92-36=", $true, $false, $false, $false, $false, $true, 1, $false, "60-56=", 2) | Out-Null
$d.Content.Find.Execute("78+9=", $true, $false, $false, $false, $false, $true, 1, $false, "73-3=", 2) | Out-Null
$d.Content.Find.Execute("62+12=", $true, $false, $false, $false, $false, $true, 1, $false, "31+2=", 2) | Out-Null
$d.Content.Find.Execute("49+39=", $true, $false, $false, $false, $false, $true, 1, $false, "62-56=", 2) | Out-Null
$d.Content.Find.Execute("87-16=", $true, $false, $false, $false, $false, $true, 1, $false, "3+71=", 2) | Out-Null
$d.Content.Find.Execute("31-30=", $true, $false, $false, $false, $false, $true, 1, $false, "21+78=", 2) | Out-Null
$d.Content.Find.Execute("10+40=", $true, $false, $false, $false, $false, $true, 1, $false, "11+7=", 2) | Out-Null
$d.Content.Find.Execute("97-49=", $true, $false, $false, $false, $false, $true, 1, $false, "63-17=", 2) | Out-Null
$d.Content.Find.Execute("8+20=", $true, $false, $false, $false, $false, $true, 1, $false, "3+24=", 2) | Out-Null
$d.Content.Find.Execute("90-63=", $true, $false, $false, $false, $false, $true, 1, $false, "92-56=", 2) | Out-Null
$d.Content.Find.Execute("25+71=", $true, $false, $false, $false, $false, $true, 1, $false, "34-13=", 2) | Out-Null
$d.Content.Find.Execute("61-34=", $true, $false, $false, $false, $false, $true, 1, $false, "46-20=", 2) | Out-Null
$d.Content.Find.Execute("99-45=", $true, $false, $false, $false, $false, $true, 1, $false, "24-18=", 2) | Out-Null
$d.Content.Find.Execute("73-18=", $true, $false, $false, $false, $false, $true, 1, $false, "7+5=", 2) | Out-Null
$d.Content.Find.Execute("85+7=", $true, $false, $false, $false, $false, $true, 1, $false, "62+22=", 2) | Out-Null
$d.Content.Find.Execute("0+25=", $true, $false, $false, $false, $false, $true, 1, $false, "45+38=", 2) | Out-Null
$d.Content.Find.Execute("67-66=", $true, $false, $false, $false, $false, $true, 1, $false, "5+5=", 2) | Out-Null
$d.Content.Find.Execute("92-12=", $true, $false, $false, $false, $false, $true, 1, $false, "13+51=", 2) | Out-Null
$d.Content.Find.Execute("88-65=", $true, $false, $false, $false, $false, $true, 1, $false, "55-39=", 2) | Out-Null
$d.Content.Find.Execute("46-22=", $true, $false, $false, $false, $false, $true, 1, $false, "37+3=", 2) | Out-Null
$d.Content.Find.Execute("50-41=", $true, $false, $false, $false, $false, $true, 1, $false, "58-17=", 2) | Out-Null
$d.Content.Find.Execute("91-42=", $true, $false, $false, $false, $false, $true, 1, $false, "85-42=", 2) | Out-Null
$d.Content.Find.Execute("21+45=", $true, $false, $false, $false, $false, $true, 1, $false, "92-86=", 2) | Out-Null
$d.Content.Find.Execute("82-47=", $true, $false, $false, $false, $false, $true, 1, $false, "19+69=", 2) | Out-Null
$d.Content.Find.Execute("6+81=", $true, $false, $false, $false, $false, $true, 1, $false, "57-19=", 2) | Out-Null
$d.Content.Find.Execute("93-91=", $true, $false, $false, $false, $false, $true, 1, $false, "72-52=", 2) | Out-Null
$d.Content.Find.Execute("63-42=", $true, $false, $false, $false, $false, $true, 1, $false, "94-90=", 2) | Out-Null
$d.Content.Find.Execute("91-30=", $true, $false, $false, $false, $false, $true, 1, $false, "45+16=", 2) | Out-Null
$d.Content.Find.Execute("14+27=", $true, $false, $false, $false, $false, $true, 1, $false, "91+5=", 2) | Out-Null
$d.Content.Find.Execute("30+0=", $true, $false, $false, $false, $false, $true, 1, $false, "49+9=", 2) | Out-Null
$d.Content.Find.Execute("44-27=", $true, $false, $false, $false, $false, $true, 1, $false, "36-30=", 2) | Out-Null
$d.Content.Find.Execute("15+73=", $true, $false, $false, $false, $false, $true, 1, $false, "51+16=", 2) | Out-Null
$d.Content.Find.Execute("77-71=", $true, $false, $false, $false, $false, $true, 1, $false, "51-35=", 2) | Out-Null
$d.Content.Find.Execute("30+63=", $true, $false, $false, $false, $false, $true, 1, $false, "54+18=", 2) | Out-Null
$d.Content.Find.Execute("90-24=", $true, $false, $false, $false, $false, $true, 1, $false, "63-56=", 2) | Out-Null
$d.Content.Find.Execute("37-26=", $true, $false, $false, $false, $false, $true, 1, $false, "67-39=", 2) | Out-Null
$d.Content.Find.Execute("88+11=", $true, $false, $false, $false, $false, $true, 1, $false, "1+65=", 2) | Out-Null
$d.Content.Find.Execute("50-5=", $true, $false, $false, $false, $false, $true, 1, $false, "75-39=", 2) | Out-Null
$d.Content.Find.Execute("50-34=", $true, $false, $false, $false, $false, $true, 1, $false, "40+2=", 2) | Out-Null
$d.Content.Find.Execute("60+4=", $true, $false, $false, $false, $false, $true, 1, $false, "75-42=", 2) | Out-Null
$d.Content.Find.Execute("59+10=", $true, $false, $false, $false, $false, $true, 1, $false, "16+25=", 2) | Out-Null
$d.Content.Find.Execute("80-35=", $true, $false, $false, $false, $false, $true, 1, $false, "31+18=", 2) | Out-Null
$d.Content.Find.Execute("80-15=", $true, $false, $false, $false, $false, $true, 1, $false, "78-51=", 2) | Out-Null
$d.Content.Find.Execute("28+35=", $true, $false, $false, $false, $false, $true, 1, $false, "18+1=", 2) | Out-Null
$d.Content.Find.Execute("58+33=", $true, $false, $false, $false, $false, $true, 1, $false, "49-24=", 2) | Out-Null
$d.Content.Find.Execute("14+21=", $true, $false, $false, $false, $false, $true, 1, $false, "35+36=", 2) | Out-Null
$d.Content.Find.Execute("66-41=", $true, $false, $false, $false, $false, $true, 1, $false, "25+26=", 2) | Out-Null
$d.Content.Find.Execute("58-16=", $true, $false, $false, $false, $false, $true, 1, $false, "76-24=", 2) | Out-Null
$d.Content.Find.Execute("63-27=", $true, $false, $false, $false, $false, $true, 1, $false, "65+27=", 2) | Out-Null
$d.Content.Find.Execute("72-11=", $true, $false, $false, $false, $false, $true, 1, $false, "26+21=", 2) | Out-Null
$d.Content.Find.Execute("43-7=", $true, $false, $false, $false, $false, $true, 1, $false, "9+81=", 2) | Out-Null
$d.Content.Find.Execute("97-32=", $true, $false, $false, $false, $false, $true, 1, $false, "49+27=", 2) | Out-Null
$d.Content.Find.Execute("59-30=", $true, $false, $false, $false, $false, $true, 1, $false, "60-54=", 2) | Out-Null
$d.Content.Find.Execute("14-14=", $true, $false, $false, $false, $false, $true, 1, $false, "82+13=", 2) | Out-Null
$d.Content.Find.Execute("76-40=", $true, $false, $false, $false, $false, $true, 1, $false, "73+17=", 2) | Out-Null
$d.Content.Find.Execute("7+73=", $true, $false, $false, $false, $false, $true, 1, $false, "84-37=", 2) | Out-Null
$d.Content.Find.Execute("67-29=", $true, $false, $false, $false, $false, $true, 1, $false, "85-30=", 2) | Out-Null
$d.Content.Find.Execute("25-12=", $true, $false, $false, $false, $false, $true, 1, $false, "53+39=", 2) | Out-Null
$d.Content.Find.Execute("59-7=", $true, $false, $false, $false, $false, $true, 1, $false, "75+6=", 2) | Out-Null
$d.Content.Find.Execute("83-15=", $true, $false, $false, $false, $false, $true, 1, $false, "93-20=", 2) | Out-Null
$d.Content.Find.Execute("93-55=", $true, $false, $false, $false, $false, $true, 1, $false, "75-33=", 2) | Out-Null
$d.Content.Find.Execute("32+48=", $true, $false, $false, $false, $false, $true, 1, $false, "60+7=", 2) | Out-Null
$d.Content.Find.Execute("37+2=", $true, $false, $false, $false, $false, $true, 1, $false, "77-26=", 2) | Out-Null
$d.Content.Find.Execute("30-2=", $true, $false, $false, $false, $false, $true, 1, $false, "5+63=", 2) | Out-Null
$d.Content.Find.Execute("7+2=", $true, $false, $false, $false, $false, $true, 1, $false, "4+71=", 2) | Out-Null
$d.Content.Find.Execute("53-47=", $true, $false, $false, $false, $false, $true, 1, $false, "37-5=", 2) | Out-Null
$d.Content.Find.Execute("33-14=", $true, $false, $false, $false, $false, $true, 1, $false, "78+17=", 2) | Out-Null
$d.Content.Find.Execute("50+48=", $true, $false, $false, $false, $false, $true, 1, $false, "53-11=", 2) | Out-Null
$d.Content.Find.Execute("78+4=", $true, $false, $false, $false, $false, $true, 1, $false, "15+56=", 2) | Out-Null
$d.Content.Find.Execute("51+28=", $true, $false, $false, $false, $false, $true, 1, $false, "34+47=", 2) | Out-Null
$d.Content.Find.Execute("65+24=", $true, $false, $false, $false, $false, $true, 1, $false, "5-0=", 2) | Out-Null
$d.Content.Find.Execute("71-41=", $true, $false, $false, $false, $false, $true, 1, $false, "57+0=", 2) | Out-Null
$d.Content.Find.Execute("69-61=", $true, $false, $false, $false, $false, $true, 1, $false, "56-26=", 2) | Out-Null
$d.Content.Find.Execute("24-3=", $true, $false, $false, $false, $false, $true, 1, $false, "64-22=", 2) | Out-Null
$d.Content.Find.Execute("17+28=", $true, $false, $false, $false, $false, $true, 1, $false, "77-56=", 2) | Out-Null
$d.Content.Find.Execute("53-45=", $true, $false, $false, $false, $false, $true, 1, $false, "26-4=", 2) | Out-Null
$d.Content.Find.Execute("54-25=", $true, $false, $false, $false, $false, $true, 1, $false, "3+39=", 2) | Out-Null
$d.Content.Find.Execute("35+0=", $true, $false, $false, $false, $false, $true, 1, $false, "12+81=", 2) | Out-Null
$d.Content.Find.Execute("35+58=", $true, $false, $false, $false, $false, $true, 1, $false, "5+57=", 2) | Out-Null
$d.Content.Find.Execute("42-38=", $true, $false, $false, $false, $false, $true, 1, $false, "55-16=", 2) | Out-Null
$d.Content.Find.Execute("26+68=", $true, $false, $false, $false, $false, $true, 1, $false, "17-12=", 2) | Out-Null
$d.Content.Find.Execute("4+57=", $true, $false, $false, $false, $false, $true, 1, $false, "54-46=", 2) | Out-Null
$d.Content.Find.Execute("26+59=", $true, $false, $false, $false, $false, $true, 1, $false, "92-66=", 2) | Out-Null
$d.Content.Find.Execute("94-7=", $true, $false, $false, $false, $false, $true, 1, $false, "88-28=", 2) | Out-Null
$d.Content.Find.Execute("27-9=", $true, $false, $false, $false, $false, $true, 1, $false, "64-12=", 2) | Out-Null
$d.Content.Find.Execute("50+15=", $true, $false, $false, $false, $false, $true, 1, $false, "54+7=", 2) | Out-Null
$d.Content.Find.Execute("58-51=", $true, $false, $false, $false, $false, $true, 1, $false, "41-8=", 2) | Out-Null
$d.Content.Find.Execute("96-3=", $true, $false, $false, $false, $false, $true, 1, $false, "28+69=", 2) | Out-Null
$d.Content.Find.Execute("9+62=", $true, $false, $false, $false, $false, $true, 1, $false, "44+43=", 2) | Out-Null
$d.Content.Find.Execute("76+23=", $true, $false, $false, $false, $false, $true, 1, $false, "20+24=", 2) | Out-Null
$d.Content.Find.Execute("99-16=", $true, $false, $false, $false, $false, $true, 1, $false, "21+17=", 2) | Out-Null
$d.Content.Find.Execute("18+8=", $true, $false, $false, $false, $false, $true, 1, $false, "46+3=", 2) | Out-Null
